$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$new = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$new.Name = "corrects.info"

# Apply the bold/bordered/centered header style (style index 1, already present in styles.xml,
# used by the "corrects" sheet header row) to the header row and the flag column.
$ws1.Range("A1").Copy()
$new.Range("B1:M1").PasteSpecial(-4122)
$new.Range("A2:A20").PasteSpecial(-4122)

# Header row (left to right)
$new.Range("B1").Value = "Algorithm"
$new.Range("C1").Value = "Comprehension"
$new.Range("D1").Value = "rec correct"
$new.Range("E1").Value = "rec total Answers"
$new.Range("F1").Value = "rec percent"
$new.Range("G1").Value = "R_ResponseTimeMean"
$new.Range("H1").Value = "R_ResponseTimeDeviation"
$new.Range("I1").Value = "it correct"
$new.Range("J1").Value = "it total Answers"
$new.Range("K1").Value = "it percent"
$new.Range("L1").Value = "I_ResponseTimeMean"
$new.Range("M1").Value = "I_ResponseTimeDeviation"

# Column B (top to bottom)
$new.Range("B2").Value = "BinarySearch"
$new.Range("B3").Value = "BinarySearch"
$new.Range("B4").Value = "BubbleSort"
$new.Range("B5").Value = "BubbleSort"
$new.Range("B6").Value = "Factorial"
$new.Range("B7").Value = "Factorial"
$new.Range("B8").Value = "Fibonacci"
$new.Range("B9").Value = "Fibonacci"
$new.Range("B10").Value = "IntegerBinary"
$new.Range("B11").Value = "IntegerBinary"
$new.Range("B12").Value = "MultiplyMatrix"
$new.Range("B13").Value = "MultiplyMatrix"
$new.Range("B14").Value = "PrimeFactors"
$new.Range("B15").Value = "PrimeFactors"
$new.Range("B16").Value = "ReverseString"
$new.Range("B17").Value = "ReverseString"
$new.Range("B18").Value = "Total"
$new.Range("B19").Value = "Total"
$new.Range("B20").Value = "Total"

# Column C (top to bottom)
$new.Range("C2").Value = "Bottom-Up"
$new.Range("C3").Value = "Top-Down"
$new.Range("C4").Value = "Bottom-Up"
$new.Range("C5").Value = "Top-Down"
$new.Range("C6").Value = "Bottom-Up"
$new.Range("C7").Value = "Top-Down"
$new.Range("C8").Value = "Bottom-Up"
$new.Range("C9").Value = "Top-Down"
$new.Range("C10").Value = "Bottom-Up"
$new.Range("C11").Value = "Top-Down"
$new.Range("C12").Value = "Bottom-Up"
$new.Range("C13").Value = "Top-Down"
$new.Range("C14").Value = "Bottom-Up"
$new.Range("C15").Value = "Top-Down"
$new.Range("C16").Value = "Bottom-Up"
$new.Range("C17").Value = "Top-Down"
$new.Range("C18").Value = "Top-Down"
$new.Range("C19").Value = "Bottom-Up"
$new.Range("C20").Value = "Total"

# Remaining numeric columns (row by row)
$new.Range("A2").Value = 0
$new.Range("D2").Value = 7
$new.Range("E2").Value = 7
$new.Range("F2").Value = 1
$new.Range("G2").Value = 441
$new.Range("H2").Value = 201
$new.Range("I2").Value = 3
$new.Range("J2").Value = 3
$new.Range("K2").Value = 1
$new.Range("L2").Value = 283
$new.Range("M2").Value = 261

$new.Range("A3").Value = 0
$new.Range("D3").Value = 3
$new.Range("E3").Value = 3
$new.Range("F3").Value = 1
$new.Range("G3").Value = 202
$new.Range("H3").Value = 99
$new.Range("I3").Value = 9
$new.Range("J3").Value = 9
$new.Range("K3").Value = 1
$new.Range("L3").Value = 322
$new.Range("M3").Value = 301

$new.Range("A4").Value = 0
$new.Range("D4").Value = 13
$new.Range("E4").Value = 13
$new.Range("F4").Value = 1
$new.Range("G4").Value = 260
$new.Range("H4").Value = 156
$new.Range("I4").Value = 8
$new.Range("J4").Value = 8
$new.Range("K4").Value = 1
$new.Range("L4").Value = 261
$new.Range("M4").Value = 218

$new.Range("A5").Value = 0
$new.Range("D5").Value = 17
$new.Range("E5").Value = 17
$new.Range("F5").Value = 1
$new.Range("G5").Value = 27
$new.Range("H5").Value = 28
$new.Range("I5").Value = 16
$new.Range("J5").Value = 16
$new.Range("K5").Value = 1
$new.Range("L5").Value = 401
$new.Range("M5").Value = 492

$new.Range("A6").Value = 0
$new.Range("D6").Value = 24
$new.Range("E6").Value = 24
$new.Range("F6").Value = 1
$new.Range("G6").Value = 132
$new.Range("H6").Value = 251
$new.Range("I6").Value = 27
$new.Range("J6").Value = 27
$new.Range("K6").Value = 1
$new.Range("L6").Value = 139
$new.Range("M6").Value = 105

$new.Range("A7").Value = 0
$new.Range("D7").Value = 25
$new.Range("E7").Value = 25
$new.Range("F7").Value = 1
$new.Range("G7").Value = 77
$new.Range("H7").Value = 35
$new.Range("I7").Value = 19
$new.Range("J7").Value = 19
$new.Range("K7").Value = 1
$new.Range("L7").Value = 109
$new.Range("M7").Value = 78

$new.Range("A8").Value = 0
$new.Range("D8").Value = 13
$new.Range("E8").Value = 13
$new.Range("F8").Value = 1
$new.Range("G8").Value = 180
$new.Range("H8").Value = 143
$new.Range("I8").Value = 18
$new.Range("J8").Value = 18
$new.Range("K8").Value = 1
$new.Range("L8").Value = 240
$new.Range("M8").Value = 144

$new.Range("A9").Value = 0
$new.Range("D9").Value = 16
$new.Range("E9").Value = 16
$new.Range("F9").Value = 1
$new.Range("G9").Value = 200
$new.Range("H9").Value = 154
$new.Range("I9").Value = 23
$new.Range("J9").Value = 23
$new.Range("K9").Value = 1
$new.Range("L9").Value = 180
$new.Range("M9").Value = 113

$new.Range("A10").Value = 0
$new.Range("D10").Value = 7
$new.Range("E10").Value = 7
$new.Range("F10").Value = 1
$new.Range("G10").Value = 420
$new.Range("H10").Value = 516
$new.Range("I10").Value = 7
$new.Range("J10").Value = 7
$new.Range("K10").Value = 1
$new.Range("L10").Value = 399
$new.Range("M10").Value = 403

$new.Range("A11").Value = 0
$new.Range("D11").Value = 12
$new.Range("E11").Value = 12
$new.Range("F11").Value = 1
$new.Range("G11").Value = 130
$new.Range("H11").Value = 92
$new.Range("I11").Value = 12
$new.Range("J11").Value = 12
$new.Range("K11").Value = 1
$new.Range("L11").Value = 210
$new.Range("M11").Value = 98

$new.Range("A12").Value = 0
$new.Range("D12").Value = 4
$new.Range("E12").Value = 4
$new.Range("F12").Value = 1
$new.Range("G12").Value = 550
$new.Range("H12").Value = 88
$new.Range("I12").Value = 8
$new.Range("J12").Value = 8
$new.Range("K12").Value = 1
$new.Range("L12").Value = 464
$new.Range("M12").Value = 357

$new.Range("A13").Value = 0
$new.Range("D13").Value = 5
$new.Range("E13").Value = 5
$new.Range("F13").Value = 1
$new.Range("G13").Value = 596
$new.Range("H13").Value = 367
$new.Range("I13").Value = 4
$new.Range("J13").Value = 4
$new.Range("K13").Value = 1
$new.Range("L13").Value = 256
$new.Range("M13").Value = 134

$new.Range("A14").Value = 0
$new.Range("D14").Value = 4
$new.Range("E14").Value = 4
$new.Range("F14").Value = 1
$new.Range("G14").Value = 142
$new.Range("H14").Value = 71
$new.Range("I14").Value = 3
$new.Range("J14").Value = 3
$new.Range("K14").Value = 1
$new.Range("L14").Value = 510
$new.Range("M14").Value = 369

$new.Range("A15").Value = 0
$new.Range("D15").Value = 5
$new.Range("E15").Value = 5
$new.Range("F15").Value = 1
$new.Range("G15").Value = 315
$new.Range("H15").Value = 371
$new.Range("I15").Value = 5
$new.Range("J15").Value = 5
$new.Range("K15").Value = 1
$new.Range("L15").Value = 166
$new.Range("M15").Value = 133

$new.Range("A16").Value = 0
$new.Range("D16").Value = 15
$new.Range("E16").Value = 15
$new.Range("F16").Value = 1
$new.Range("G16").Value = 217
$new.Range("H16").Value = 190
$new.Range("I16").Value = 22
$new.Range("J16").Value = 22
$new.Range("K16").Value = 1
$new.Range("L16").Value = 224
$new.Range("M16").Value = 237

$new.Range("A17").Value = 0
$new.Range("D17").Value = 12
$new.Range("E17").Value = 12
$new.Range("F17").Value = 1
$new.Range("G17").Value = 250
$new.Range("H17").Value = 264
$new.Range("I17").Value = 23
$new.Range("J17").Value = 23
$new.Range("K17").Value = 1
$new.Range("L17").Value = 160
$new.Range("M17").Value = 107

$new.Range("A18").Value = 0
$new.Range("D18").Value = 87
$new.Range("E18").Value = 87
$new.Range("F18").Value = 1
$new.Range("G18").Value = 241
$new.Range("H18").Value = 256
$new.Range("I18").Value = 96
$new.Range("J18").Value = 96
$new.Range("K18").Value = 1
$new.Range("L18").Value = 250
$new.Range("M18").Value = 241

$new.Range("A19").Value = 0
$new.Range("D19").Value = 95
$new.Range("E19").Value = 95
$new.Range("F19").Value = 1
$new.Range("G19").Value = 161
$new.Range("H19").Value = 207
$new.Range("I19").Value = 111
$new.Range("J19").Value = 111
$new.Range("K19").Value = 1
$new.Range("L19").Value = 212
$new.Range("M19").Value = 238

$new.Range("A20").Value = 0
$new.Range("D20").Value = 182
$new.Range("E20").Value = 182
$new.Range("F20").Value = 1
$new.Range("G20").Value = 199
$new.Range("H20").Value = 235
$new.Range("I20").Value = 207
$new.Range("J20").Value = 207
$new.Range("K20").Value = 1
$new.Range("L20").Value = 230
$new.Range("M20").Value = 239

# Keep the original "corrects" sheet as the active/selected tab, matching the source workbook
$ws1.Activate()
